$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Zero" config
$ws.Range("A4").Value = "Zero"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3

# Row 5: "ZeroAcc" config
$ws.Range("A5").Value = "ZeroAcc"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 3

$ws.Range("C6").Select()
